$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B24: was stored as text "2", should be a real number 2
$ws.Cells.Item(24, 2).Value = 2

# Add new row 25 with annotation data
$ws.Cells.Item(25, 1).Value = "Sunsi Wu"

# B25 must stay as text "3" (not a number), so force text format first
$b25 = $ws.Cells.Item(25, 2)
$b25.NumberFormat = "@"
$b25.Value = "3"
$b25.ClearFormats()

$ws.Cells.Item(25, 3).Value = "无"
$ws.Cells.Item(25, 4).Value = "FBK"
$ws.Cells.Item(25, 5).Value = "OTH"
$ws.Cells.Item(25, 6).Value = "91b1b71f-4957-400a-bdb5-bced2ed448de"
$ws.Cells.Item(25, 7).Value = "S1CChZ-CZ_annotated.xlsx"
$ws.Cells.Item(25, 8).Value = "It took us as a few weeks to reply because we took the time to implement as much as possible of the feedback."
